# Insert a new row at position 42; this shifts existing rows 42..86 down to 43..87
# and preserves cell styles (e.g. the date-formatted column D).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(42).Insert()

# Populate the newly inserted row 42 with the new record's data.
$ws.Range("A42").Value = 10
$ws.Range("B42").Value = "Vega Modelo de Temuco"
$ws.Range("C42").Value = "La Araucanía"
$ws.Range("D42").Value = 44484
$ws.Range("E42").Value = 9
$ws.Range("F42").Value = 100112012
$ws.Range("G42").Value = "Espinaca"
$ws.Range("H42").Value = "Sin especificar"
$ws.Range("I42").Value = "Primera"
$ws.Range("J42").Value = 20
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 10000
$ws.Range("M42").Value = 10000
$ws.Range("N42").Value = "$/docena de atados"
$ws.Range("O42").Value = "Región de La Araucanía"
$ws.Range("P42").Value = 3333
$ws.Range("Q42").Value = 3
$ws.Range("R42").Value = "Hortaliza"
